# ROM.xlsx keypoints export — add more joints (left/right hip & knee split
# out, shoulder flexion added) and drop the old scratch rows, leaving a
# single two-row table: a header row of joint names and one row of values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old sheet had header text in row 1 and leftover sample/placeholder
# data in rows 2-5 (columns A:D). Wipe all of it before laying out the new,
# wider table.
$ws.Range("A1:D5").ClearContents()

# New header row: Trunk / Right+Left Hip / Right+Left Knee / Ankle /
# Right+Left Shoulder flexion, one column each (A1:H1).
$ws.Range("A1").Value = "TrunkROM"
$ws.Range("B1").Value = "RIGHT HIP ROM"
$ws.Range("C1").Value = "LEFT HIP ROM"
$ws.Range("D1").Value = "RIGHT KNEE ROM"
$ws.Range("E1").Value = "LEFT KNEE ROM"
$ws.Range("F1").Value = "AnkleROM"
$ws.Range("G1").Value = "RIGHT SHOULDER FLEXION ROM"
$ws.Range("H1").Value = "LEFT SHOULDER FLEXION ROM"

# Give the four new header cells (E1:H1) the same bold/centered/bordered
# look already used for the original header cells (A1:D1).
$ws.Range("A1:D1").Copy()
$ws.Range("E1:H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 2: the measured ROM values for this frame (right hip, right knee and
# right shoulder flexion are populated; the rest are left blank).
$ws.Range("B2").Value = 0.269439697265625
$ws.Range("D2").Value = 0.0922698974609375
$ws.Range("G2").Value = 0.9779548645019531
